$wb = $excel.ActiveWorkbook

# Work on the "Sucursales" sheet (sheet2)
$ws = $wb.Worksheets.Item("Sucursales")

# Add a new header in C1: "Id zona", using the same bold/centered header style as A1/B1
$ws.Range("C1").Value = "Id zona"
$ws.Range("C1").Style = $ws.Range("A1").Style
$ws.Cells.Item(1, 3).HorizontalAlignment = -4108   # xlCenter
$ws.Cells.Item(1, 3).VerticalAlignment = -4108     # xlCenter

# Widen column C and center its contents like the existing style used elsewhere
$ws.Columns.Item(3).ColumnWidth = 11.43
$ws.Columns.Item(3).HorizontalAlignment = -4108
$ws.Columns.Item(3).VerticalAlignment = -4108

# Move/keep the active selection on this sheet at E6
$ws.Activate()
$ws.Range("E6").Select()
